# Updates the cryptocurrency price list (columns D = Price, E = Volume(1h))
# to reflect the latest scraped values, as produced by the GitHub Actions
# scheduled job that refreshes cryptos.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values like "43.278.96" or "0.649" that are
# really text (thousands are dot-separated, so plain decimals such as
# "0.649" would otherwise be auto-parsed into numbers by Excel). Force the
# column to Text before writing, then restore its original (General) style
# so no stray number formatting is left behind on the cells.
$ws.Range("D2:D51").NumberFormat = "@"

$sub3 = [char]0x2083

$updates = @(
    @{ Row = 2;  D = "43.278.96";       E = "  -1.72%  " },
    @{ Row = 3;  D = "2.346.24";        E = "  +3.71%  " },
    @{ Row = 4;  D = $null;             E = "  +0.13%  " },
    @{ Row = 5;  D = "0.649";           E = "  +1.70%  " },
    @{ Row = 6;  D = "230.86";          E = "  +0.13%  " },
    @{ Row = 7;  D = "65.24";           E = "  +1.33%  " },
    @{ Row = 8;  D = $null;             E = "  +0.11%  " },
    @{ Row = 9;  D = "0.456";           E = "  +1.70%  " },
    @{ Row = 10; D = "0.0943";          E = "  -5.31%  " },
    @{ Row = 11; D = "56.90";           E = "  -0.61%  " },
    @{ Row = 12; D = "26.58";           E = "  -1.52%  " },
    @{ Row = 13; D = "2.694.17";        E = "  +3.58%  " },
    @{ Row = 14; D = "0.104";           E = "  -1.47%  " },
    @{ Row = 15; D = "15.28";           E = "  -2.53%  " },
    @{ Row = 16; D = "6.23";            E = "  +2.28%  " },
    @{ Row = 17; D = $null;             E = "  +0.27%  " },
    @{ Row = 18; D = "2.344.84";        E = "  +3.75%  " },
    @{ Row = 19; D = "43.248.33";       E = "  -1.42%  " },
    @{ Row = 20; D = "0.0${sub3}0971";  E = "  -3.01%  " },
    @{ Row = 21; D = "73.50";           E = "  -0.09%  " },
    @{ Row = 22; D = "6.16";            E = "  +1.35%  " },
    @{ Row = 23; D = "247.44";          E = "  -1.55%  " },
    @{ Row = 24; D = $null;             E = "  +21.01%  " },
    @{ Row = 25; D = $null;             E = "  +0.06%  " },
    @{ Row = 26; D = "2.44";            E = "  -0.41%  " },
    @{ Row = 27; D = "2.26";            E = "  +0.78%  " },
    @{ Row = 28; D = "9.85";            E = "  -2.53%  " },
    @{ Row = 29; D = "175.10";          E = "  +2.38%  " },
    @{ Row = 30; D = "22.23";           E = "  +6.30%  " },
    @{ Row = 31; D = "1.50";            E = "  +7.90%  " },
    @{ Row = 32; D = $null;             E = "  -7.97%  " },
    @{ Row = 33; D = "0.125";           E = "  +0.31%  " },
    @{ Row = 34; D = "4.96";            E = "  +3.55%  " },
    @{ Row = 35; D = "0.0686";          E = "  -2.97%  " },
    @{ Row = 36; D = "4.97";            E = "  +1.31%  " },
    @{ Row = 37; D = "2.48";            E = "  +7.34%  " },
    @{ Row = 38; D = $null;             E = "  -0.66%  " },
    @{ Row = 39; D = $null;             E = "  -5.81%  " },
    @{ Row = 41; D = $null;             E = "  +0.08%  " },
    @{ Row = 42; D = $null;             E = "  +7.93%  " },
    @{ Row = 43; D = "17.83";           E = "  +2.86%  " },
    @{ Row = 44; D = $null;             E = "  +7.16%  " },
    @{ Row = 45; D = "98.32";           E = "  +0.08%  " },
    @{ Row = 46; D = $null;             E = "  -0.58%  " },
    @{ Row = 47; D = $null;             E = "  -1.54%  " },
    @{ Row = 48; D = "0.0943";          E = "  -3.36%  " },
    @{ Row = 49; D = "1.433.60";        E = "  -0.63%  " },
    @{ Row = 50; D = "2.567.02";        E = "  +3.68%  " },
    @{ Row = 51; D = "0.000202";        E = "  -9.41%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# Restore the column's original (unformatted) style now that all the text
# values are safely stored, so only the cell contents differ from before.
$ws.Range("D2:D51").Style = "Normal"
